# Binomiales_Wins.xlsx - add "Low" (L_Periodo1-8) and "Pre_Periodo1-4" tracking
# columns for the sessions that have been run so far, and highlight the
# already-filled "Low"/preguntas cells in yellow. Also collapse/hide the
# columns that are no longer the focus of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New header row labels (AQ1:BB1)
# ---------------------------------------------------------------------
$headers = @(
    "Pre_Periodo1","Pre_Periodo2","Pre_Periodo3","Pre_Periodo4",
    "L_Periodo1","L_Periodo2","L_Periodo3","L_Periodo4",
    "L_Periodo5","L_Periodo6","L_Periodo7","L_Periodo8"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, 43 + $i).Value = $headers[$i]
}

# ---------------------------------------------------------------------
# 2. New data for rows 2-11 (sessions 1, 2 and the first two periods of
#    session 3), columns AQ:BB
# ---------------------------------------------------------------------
$data = @(
    @(0,1,1,1,0,0,0,1,1,0,1,1),
    @(1,1,0,1,0,1,0,0,1,0,0,1),
    @(0,0,1,0,0,0,1,0,0,0,1,0),
    @(1,1,0,1,1,0,0,1,1,1,0,0),
    @(1,0,0,1,1,1,1,1,0,0,1,1),
    @(1,1,0,1,1,1,1,1,0,0,0,0),
    @(1,1,0,0,1,0,0,1,0,1,0,0),
    @(0,1,1,0,0,1,1,1,0,1,1,0),
    @(1,0,1,1,0,0,1,1,0,0,0,1),
    @(0,1,1,1,0,1,0,1,0,1,0,0)
)
for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item(2 + $r, 43 + $c).Value = $rowVals[$c]
    }
}

# ---------------------------------------------------------------------
# 3. Highlight (yellow) the cells that already contain confirmed data for
#    completed sessions: E2:E29, AB2:AB33, AC2:AC37
# ---------------------------------------------------------------------
$ws.Range("E2:E29").Interior.Color = 65535
$ws.Range("AB2:AB33").Interior.Color = 65535
$ws.Range("AC2:AC37").Interior.Color = 65535

# ---------------------------------------------------------------------
# 4. Hide the columns that are not the current focus (F:AG and AI:AP),
#    keep AQ:AT (the new Pre_Periodo columns) visible and sized.
# ---------------------------------------------------------------------
$ws.Range("F1:AG1").EntireColumn.Hidden = $true
$ws.Range("AI1:AP1").EntireColumn.Hidden = $true
$ws.Range("AI1:AP1").ColumnWidth = 10.6
$ws.Range("AQ1:AT1").ColumnWidth = 10.6

# ---------------------------------------------------------------------
# 5. Move the active selection, as it was left after the edit
# ---------------------------------------------------------------------
$ws.Range("AT17").Select()
